$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value while forcing/preserving Text storage for
# numeric-looking strings (e.g. "1.002"), which Excel would otherwise
# auto-convert to a Number. Resetting the Style back to "Normal" after
# the write drops the temporary "@" text format so no stray cell style
# is left behind.
function Set-TextValue($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

$ws.Range("D2").Value = '27.365.00'
$ws.Range("E2").Value = '  -0.85%  '

$ws.Range("D3").Value = '1.785.67'
$ws.Range("E3").Value = '  -2.17%  '

Set-TextValue "D4" '1.002'
$ws.Range("E4").Value = '  -0.05%  '

Set-TextValue "D5" '340.53'
$ws.Range("E5").Value = '  -0.28%  '

$ws.Range("E6").Value = '  +0.01%  '

Set-TextValue "D7" '0.3957'
$ws.Range("E7").Value = '  +3.37%  '

Set-TextValue "D8" '0.3470'
$ws.Range("E8").Value = '  -1.87%  '

Set-TextValue "D9" '48.02'
$ws.Range("E9").Value = '  -3.46%  '

Set-TextValue "D10" '1.196'
$ws.Range("E10").Value = '  -3.60%  '

Set-TextValue "D11" '0.07477'
$ws.Range("E11").Value = '  -3.49%  '

Set-TextValue "D12" '0.9988'
$ws.Range("E12").Value = '  -0.21%  '

Set-TextValue "D13" '21.75'
$ws.Range("E13").Value = '  -2.96%  '

Set-TextValue "D14" '6.470'
$ws.Range("E14").Value = '  -2.38%  '

$ws.Range("D15").Value = '1.783.30'
$ws.Range("E15").Value = '  -2.14%  '

Set-TextValue "D16" '7.109'
$ws.Range("E16").Value = '  -1.36%  '

Set-TextValue "D17" '0.00001096'
$ws.Range("E17").Value = '  -2.85%  '

Set-TextValue "D18" '0.06692'
$ws.Range("E18").Value = '  -0.47%  '

Set-TextValue "D19" '84.70'
$ws.Range("E19").Value = '  -2.98%  '

Set-TextValue "D20" '0.9992'
$ws.Range("E20").Value = '  -0.15%  '

Set-TextValue "D21" '17.76'
$ws.Range("E21").Value = '  +0.34%  '

Set-TextValue "D22" '6.505'
$ws.Range("E22").Value = '  -0.69%  '

$ws.Range("D23").Value = '27.347.81'
$ws.Range("E23").Value = '  -0.89%  '

Set-TextValue "D24" '12.39'
$ws.Range("E24").Value = '  -5.90%  '

Set-TextValue "D25" '2.393'
$ws.Range("E25").Value = '  -3.42%  '

Set-TextValue "D26" '21.23'
$ws.Range("E26").Value = '  -4.16%  '

$ws.Range("B27").Value = 'ImmutableX'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue "D27" '1.461'
$ws.Range("E27").Value = '  -1.55%  '

$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue "D28" '2.488'
$ws.Range("E28").Value = '  -7.19%  '

Set-TextValue "D29" '157.51'
$ws.Range("E29").Value = '  +2.91%  '

$ws.Range("D30").Value = '1.984.43'
$ws.Range("E30").Value = '  -2.23%  '

Set-TextValue "D31" '135.96'
$ws.Range("E31").Value = '  +0.19%  '

Set-TextValue "D32" '4.028'
$ws.Range("E32").Value = '  -1.40%  '

Set-TextValue "D33" '5.970'
$ws.Range("E33").Value = '  -6.02%  '

Set-TextValue "D34" '0.08829'
$ws.Range("E34").Value = '  +0.21%  '

Set-TextValue "D35" '12.98'
$ws.Range("E35").Value = '  -7.12%  '

Set-TextValue "D36" '0.02424'
$ws.Range("E36").Value = '  +0.51%  '

$ws.Range("E37").Value = '  -4.67%  '

Set-TextValue "D38" '5.408'
$ws.Range("E38").Value = '  -4.38%  '

Set-TextValue "D39" '0.06478'
$ws.Range("E39").Value = '  -0.47%  '

Set-TextValue "D40" '0.6826'
$ws.Range("E40").Value = '  -3.12%  '

Set-TextValue "D41" '0.2209'
$ws.Range("E41").Value = '  -2.57%  '

Set-TextValue "D42" '1.254'
$ws.Range("E42").Value = '  -3.29%  '

Set-TextValue "D43" '8.339'
$ws.Range("E43").Value = '  -8.99%  '

Set-TextValue "D44" '14.48'
$ws.Range("E44").Value = '  -1.79%  '

Set-TextValue "D45" '0.9989'
$ws.Range("E45").Value = '  -0.12%  '

Set-TextValue "D46" '0.6380'
$ws.Range("E46").Value = '  -3.78%  '

Set-TextValue "D47" '3.878'
$ws.Range("E47").Value = '  -1.37%  '

Set-TextValue "D48" '2.135'
$ws.Range("E48").Value = '  -2.49%  '

Set-TextValue "D49" '132.05'
$ws.Range("E49").Value = '  -1.09%  '

Set-TextValue "D50" '0.07151'
$ws.Range("E50").Value = '  -2.40%  '

Set-TextValue "D51" '79.18'
$ws.Range("E51").Value = '  -2.81%  '
